# B6-PowerPoint.pptx edit
# 1) Three tables (slides 14, 15, 16) get their table style switched from the
#    custom "Table_0" style ({19C7FDCC-C8C6-4910-B60A-49EA0BDFDA7B}) to the
#    built-in style {458D5E36-5DB5-4319-BF20-1EF5BDFEA904} (as chosen from the
#    Table Design gallery in PowerPoint).
# 2) The deck's theme switches from the "Integral" (Red Violet) look back to
#    the default "Office Theme" colors - done here by recoloring the theme's
#    12 color-scheme slots to the stock Office palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -------------------------------------------------
$newTableStyleId = "{458D5E36-5DB5-4319-BF20-1EF5BDFEA904}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newTableStyleId)
    }
}

# --- 2) Switch the presentation's theme colors back to the Office defaults --------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

# Index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#              8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# (RGB values are packed as 0x00BBGGRR, matching the COM RGB() long layout)
$colorScheme.Item(1).RGB  = 0x000000    # dk1     000000
$colorScheme.Item(2).RGB  = 0xFFFFFF    # lt1     FFFFFF
$colorScheme.Item(3).RGB  = 0x6A5444    # dk2     44546A
$colorScheme.Item(4).RGB  = 0xE6E6E7    # lt2     E7E6E6
$colorScheme.Item(5).RGB  = 0xD59B5B    # accent1 5B9BD5
$colorScheme.Item(6).RGB  = 0x317DED    # accent2 ED7D31
$colorScheme.Item(7).RGB  = 0xA5A5A5    # accent3 A5A5A5
$colorScheme.Item(8).RGB  = 0x00C0FF    # accent4 FFC000
$colorScheme.Item(9).RGB  = 0xC47244    # accent5 4472C4
$colorScheme.Item(10).RGB = 0x47AD70    # accent6 70AD47
$colorScheme.Item(11).RGB = 0xC16305    # hlink   0563C1
$colorScheme.Item(12).RGB = 0x724F95    # folHlink 954F72
